$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "27.628.95"
Set-TextValue "E2" "  -2.45%  "
Set-TextValue "D3" "1.761.23"
Set-TextValue "E3" "  -3.27%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.11%  "
Set-TextValue "D5" "324.30"
Set-TextValue "E5" "  -1.40%  "
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  -0.14%  "
Set-TextValue "D7" "0.4276"
Set-TextValue "E7" "  -1.72%  "
Set-TextValue "D8" "0.3606"
Set-TextValue "E8" "  -2.21%  "
Set-TextValue "D9" "0.07559"
Set-TextValue "E9" "  -2.08%  "
Set-TextValue "D10" "42.33"
Set-TextValue "E10" "  -6.02%  "
Set-TextValue "D11" "1.109"
Set-TextValue "E11" "  -2.73%  "
Set-TextValue "D12" "1.000"
Set-TextValue "E12" "  -0.19%  "
Set-TextValue "D13" "20.84"
Set-TextValue "E13" "  -6.17%  "
Set-TextValue "D14" "6.077"
Set-TextValue "E14" "  -4.10%  "
Set-TextValue "D15" "7.242"
Set-TextValue "E15" "  -4.28%  "
Set-TextValue "D16" "1.760.47"
Set-TextValue "E16" "  -4.24%  "
Set-TextValue "D17" "93.34"
Set-TextValue "E17" "  +0.16%  "
Set-TextValue "E18" "  -1.46%  "
Set-TextValue "D19" "0.06428"
Set-TextValue "E19" "  -1.52%  "
Set-TextValue "D20" "0.9995"
Set-TextValue "E20" "  -0.18%  "
Set-TextValue "D21" "17.16"
Set-TextValue "E21" "  -2.26%  "
Set-TextValue "D22" "5.901"
Set-TextValue "E22" "  -6.31%  "
Set-TextValue "D23" "27.678.72"
Set-TextValue "E23" "  -2.44%  "
Set-TextValue "D24" "11.30"
Set-TextValue "E24" "  -3.42%  "
Set-TextValue "D25" "2.127"
Set-TextValue "E25" "  +9.28%  "
Set-TextValue "D26" "162.58"
Set-TextValue "E26" "  +0.58%  "
Set-TextValue "D27" "20.36"
Set-TextValue "E27" "  -2.51%  "
Set-TextValue "D28" "1.962.25"
Set-TextValue "E28" "  -3.80%  "
Set-TextValue "D29" "2.164"
Set-TextValue "E29" "  -6.16%  "
Set-TextValue "D30" "125.51"
Set-TextValue "E30" "  -2.97%  "
Set-TextValue "D31" "1.109"
Set-TextValue "E31" "  -9.73%  "
Set-TextValue "D32" "5.606"
Set-TextValue "E32" "  -7.30%  "
Set-TextValue "D33" "3.658"
Set-TextValue "E33" "  +2.42%  "
Set-TextValue "D34" "0.08950"
Set-TextValue "E34" "  -2.96%  "
Set-TextValue "D35" "12.24"
Set-TextValue "E35" "  -5.66%  "
Set-TextValue "D36" "0.02296"
Set-TextValue "E36" "  -3.03%  "
Set-TextValue "D38" "0.06035"
Set-TextValue "E38" "  -2.89%  "
Set-TextValue "D39" "0.6364"
Set-TextValue "E39" "  -3.77%  "
Set-TextValue "D40" "4.965"
Set-TextValue "E40" "  -5.27%  "
Set-TextValue "E41" "  -0.46%  "
Set-TextValue "D42" "0.9994"
Set-TextValue "E42" "  -0.14%  "
Set-TextValue "D43" "1.396"
Set-TextValue "E43" "  -3.00%  "
Set-TextValue "D44" "7.909"
Set-TextValue "E44" "  -3.35%  "
Set-TextValue "D45" "13.44"
Set-TextValue "E45" "  -3.70%  "
Set-TextValue "D46" "0.5953"
Set-TextValue "E46" "  -3.12%  "
Set-TextValue "D47" "3.714"
Set-TextValue "E47" "  -1.30%  "
Set-TextValue "D48" "1.996"
Set-TextValue "E48" "  -1.87%  "
Set-TextValue "D49" "123.34"
Set-TextValue "E49" "  -2.57%  "
Set-TextValue "D50" "1.171"
Set-TextValue "E50" "  +0.82%  "
Set-TextValue "D51" "0.06862"
Set-TextValue "E51" "  -2.29%  "
